$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: set text, then copy formatting (border/bold/alignment) from an existing header cell
$ws.Range("T1").Value = "WARDS_POSEES_MOYENNE"
$ws.Range("U1").Value = "WARDS_DETRUITES_MOYENNE"
$ws.Range("V1").Value = "WARDS_PINKS_MOYENNE"
$ws.Range("S1").Copy()
$ws.Range("T1:V1").PasteSpecial(-4122)

# Row 2
$ws.Range("G2").Value = 2364
$ws.Range("H2").Value = 830
$ws.Range("R2").Value = 29.92
$ws.Range("T2").Value = 10.51
$ws.Range("U2").Value = 4.9
$ws.Range("V2").Value = 6.68

# Row 3
$ws.Range("G3").Value = 889
$ws.Range("H3").Value = 442
$ws.Range("R3").Value = 21.68
$ws.Range("T3").Value = 10.78
$ws.Range("U3").Value = 2.05
$ws.Range("V3").Value = 2.27

# Row 4
$ws.Range("D4").Value = 297
$ws.Range("F4").Value = 8648.3755
$ws.Range("G4").Value = 5961
$ws.Range("H4").Value = 2559
$ws.Range("I4").Value = 872
$ws.Range("J4").Value = 798
$ws.Range("K4").Value = 58978
$ws.Range("L4").Value = 2098
$ws.Range("M4").Value = 1645
$ws.Range("N4").Value = 2267
$ws.Range("O4").Value = 7.063973063973064
$ws.Range("P4").Value = 5.538720538720539
$ws.Range("Q4").Value = 7.632996632996633
$ws.Range("R4").Value = 20.07
$ws.Range("S4").Value = 1747.15
$ws.Range("T4").Value = 8.619999999999999
$ws.Range("U4").Value = 2.94
$ws.Range("V4").Value = 2.69

# Row 5
$ws.Range("G5").Value = 5572
$ws.Range("H5").Value = 2575
$ws.Range("R5").Value = 42.53
$ws.Range("T5").Value = 19.66
$ws.Range("U5").Value = 5.15
$ws.Range("V5").Value = 5.59

# Row 6
$ws.Range("G6").Value = 2811
$ws.Range("H6").Value = 1130
$ws.Range("R6").Value = 57.37
$ws.Range("T6").Value = 23.06
$ws.Range("U6").Value = 7.59
$ws.Range("V6").Value = 7.22

# Row 7
$ws.Range("G7").Value = 553
$ws.Range("R7").Value = 26.33
$ws.Range("T7").Value = 12.1
$ws.Range("U7").Value = 2.14
$ws.Range("V7").Value = 2.71

# Row 8
$ws.Range("G8").Value = 4593
$ws.Range("H8").Value = 2191
$ws.Range("R8").Value = 80.58
$ws.Range("T8").Value = 38.44
$ws.Range("U8").Value = 11.98
$ws.Range("V8").Value = 7.7

# Row 9
$ws.Range("G9").Value = 1895
$ws.Range("H9").Value = 641
$ws.Range("R9").Value = 35.09
$ws.Range("T9").Value = 11.87
$ws.Range("U9").Value = 4.94
$ws.Range("V9").Value = 2.93

# Row 10
$ws.Range("G10").Value = 7693
$ws.Range("H10").Value = 4159
$ws.Range("R10").Value = 63.58
$ws.Range("T10").Value = 34.37
$ws.Range("U10").Value = 8.15
$ws.Range("V10").Value = 8.279999999999999

# Row 11
$ws.Range("G11").Value = 1623
$ws.Range("H11").Value = 685
$ws.Range("R11").Value = 33.12
$ws.Range("T11").Value = 13.98
$ws.Range("U11").Value = 2.88
$ws.Range("V11").Value = 5.12

# Row 12
$ws.Range("G12").Value = 277
$ws.Range("H12").Value = 120
$ws.Range("R12").Value = 21.31
$ws.Range("T12").Value = 9.23
$ws.Range("U12").Value = 1.69
$ws.Range("V12").Value = 0.31

# Row 13
$ws.Range("G13").Value = 720
$ws.Range("R13").Value = 26.67
$ws.Range("T13").Value = 12.44
$ws.Range("U13").Value = 3.19
$ws.Range("V13").Value = 3.96
